$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the filter/header column from "no_of_arrear" to "history_of_arrear"
$ws.Range("I1").Value = "history_of_arrear"

# Update the active selection to reflect the new working cell
$ws.Range("I2").Select()
